$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2; this pushes the existing data rows (2-28)
# down to (3-29) and keeps the header in row 1.
$ws.Range("A2").EntireRow.Insert()

# The inserted row picks up formatting copied from the row above (the bold
# header) - strip that so it matches the plain body-row formatting used
# throughout the rest of the table.
$ws.Range("A2:T2").ClearFormats()

# New weekly record for "Feria Lagunitas de Puerto Montt" (most of the
# descriptive columns repeat the values already used by the record that
# used to sit in row 2, which is now row 3).
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C2").Value = "Los Lagos"

$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").Value = 44901

$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100101
$ws.Range("H2").Value = "Berries"
$ws.Range("I2").Value = 100101001
$ws.Range("J2").Value = "Arándano (blue)"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 400
$ws.Range("N2").Value = 5000
$ws.Range("O2").Value = 5500
$ws.Range("P2").Value = 5250
$ws.Range("Q2").Value = "$/bandeja 2 kilos"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 2625
$ws.Range("T2").Value = 2
